# KIBON-120: Gesuch-Stichtag statistic translated
#
# The "Data" sheet header row used to contain hard-coded German labels.
# These are replaced with placeholder tokens that the report engine fills
# in with the translated column titles at render time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A1").Value = "{bgNummerTitle}"
$ws.Range("B1").Value = "{institutionTitle}"
$ws.Range("C1").Value = "{angebotTitle}"
$ws.Range("D1").Value = "{periodeTitle}"
$ws.Range("E1").Value = "{gesuchLaufNrTitle}"
$ws.Range("F1").Value = "{nichtFreigegebenTitle}"
$ws.Range("G1").Value = "{mahnungenTitle}"
$ws.Range("H1").Value = "{beschwerdeTitle}"

# The active selection on the "Data" sheet moved to H2.
[void]$ws.Range("H2").Select()
